$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Q0) ---
$ws.Range("B2").Value2 = -0.1380044728926245
$ws.Range("C2").Value2 = 0.4405518369217296
$ws.Range("D2").Value2 = 0.3098877474931226
$ws.Range("E2").Value2 = 0.5566756214287838
$ws.Range("F2").Value2 = 0.5596561568957058

# --- Row 3 (Q1) ---
$ws.Range("B3").Value2 = -0.05072747431594977
$ws.Range("C3").Value2 = 0.3686789333054596
$ws.Range("D3").Value2 = 0.2339459977537088
$ws.Range("E3").Value2 = 0.483679643724758
$ws.Range("F3").Value2 = 0.5006533543232978

# --- Row 4 (Q2) ---
$ws.Range("B4").Value2 = -0.05721685844930368
$ws.Range("C4").Value2 = 0.4206015753809266
$ws.Range("D4").Value2 = 0.2672981648104081
$ws.Range("E4").Value2 = 0.5170088633770297
$ws.Range("F4").Value2 = 0.5366811099065004

# --- Row 5 (Q3) ---
$ws.Range("B5").Value2 = -0.07739811459342891
$ws.Range("C5").Value2 = 0.4504919024182757
$ws.Range("D5").Value2 = 0.251819281882534
$ws.Range("E5").Value2 = 0.5018159840843395
$ws.Range("F5").Value2 = 0.5200112451802442

# --- Row 6 (Q4) ---
$ws.Range("B6").Value2 = -0.05962591978527489
$ws.Range("C6").Value2 = 0.4241150868043847
$ws.Range("D6").Value2 = 0.2659677999092564
$ws.Range("E6").Value2 = 0.515720660735302
$ws.Range("F6").Value2 = 0.5399717581082021

# --- Row 7 (Q5) ---
$ws.Range("B7").Value2 = -0.06396073784541235
$ws.Range("C7").Value2 = 0.4341613455797451
$ws.Range("D7").Value2 = 0.2840886535732064
$ws.Range("E7").Value2 = 0.5329996750216706
$ws.Range("F7").Value2 = 0.5612462804205578
$ws.Range("G7").Value2 = 9

# --- Row 8 (Q6) ---
$ws.Range("B8").Value2 = 0.01200689806893836
$ws.Range("C8").Value2 = 0.4563415053098694
$ws.Range("D8").Value2 = 0.2743301623456585
$ws.Range("E8").Value2 = 0.5237653695555468
$ws.Range("F8").Value2 = 0.5736054358993686
$ws.Range("G8").Value2 = 6

# --- Row 9 (Q7) ---
$ws.Range("B9").Value2 = -0.5577316861920555
$ws.Range("C9").Value2 = 0.5577316861920555
$ws.Range("D9").Value2 = 0.3432890341754003
$ws.Range("E9").Value2 = 0.5859087251231204
$ws.Range("F9").Value2 = 0.2198558632130385
$ws.Range("G9").Value2 = 3

# --- Row 10 (new row, Q8) ---
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value2 = -0.7745672082728081
$ws.Range("C10").Value2 = 0.7745672082728081
$ws.Range("D10").Value2 = 0.5999543601315317
$ws.Range("E10").Value2 = 0.7745672082728081
$ws.Range("G10").Value2 = 1
